$d = $word.ActiveDocument

$replacements = @(
    @("41×30=", "43×90="),
    @("77×99=", "43×14="),
    @("94×15=", "80×56="),
    @("19×60=", "68×65="),
    @("98×75=", "99×21="),
    @("71×85=", "27×51="),
    @("45×81=", "41×64="),
    @("62×78=", "72×57="),
    @("48×52=", "74×16="),
    @("50×64=", "53×39="),
    @("51×31=", "65×93="),
    @("34×53=", "22×70="),
    @("83×15=", "33×46="),
    @("81×52=", "18×46="),
    @("78×72=", "43×82="),
    @("16×98=", "14×17="),
    @("14×58=", "55×20="),
    @("80×30=", "79×79="),
    @("73×35=", "27×11="),
    @("75×32=", "61×74="),
    @("98×57=", "44×53="),
    @("99×46=", "48×68="),
    @("52×59=", "59×73="),
    @("74×48=", "59×75="),
    @("32×32=", "89×68=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
